# Merge the two separate "battery temps left/right" CAN-IDs (rows 6 & 7)
# into a single "CAN_ID_BATTERY_TEMPS" row, as part of adding the JetiBox
# config to the protocol sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 becomes the merged "battery temps" entry (keeps the 0xE0 id,
# drops the _LEFT/_RIGHT split).
$ws.Range("A6").Value = "CAN_ID_BATTERY_TEMPS"
$ws.Range("B6").Value = "0xE0"
$ws.Range("C6").Value = "Akku-Temperaturen"
$ws.Range("D6").Value = 8

# The old row 7 (CAN_ID_BATTERY_TEMPS_RIGHT / 0xE1) is no longer needed;
# remove it and shift the rest of the sheet up.
$ws.Rows("7").Delete()

# Match the author's last-saved selection.
$ws.Range("C7").Select()
